# Bulk-insert the remaining machine types below the existing list
# (A1:A6 already contains the header + first five entries).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Vertical tube bending",
    "Shearing",
    "Sheet bending",
    "Notching machine",
    "Mechanical press machine",
    "hydraullic press machine",
    "Speedy seamer",
    "Drilling machine",
    "Vertical band saw"
)

$row = 7
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}

# Widen column A so the new, longer entries stay fully visible.
$ws.Columns.Item(1).ColumnWidth = 22.33203125

# Leave the selection on the second-to-last inserted cell, as in the source edit.
$ws.Range("A14").Select() | Out-Null
